$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSI5RB6ZGB%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T190011Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJIMEYCIQCw9VcB8s63bA3LOPNHjnRYlPuYo07C22Nhbw7CzBsyogIhAMo53igAGwPughrSLjn9PT2aC%2Bn8cieuF2ARYHGRU0Q%2BKuYDCJH%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgyT3Etg0ZfLnQA5s%2BAqugPb9%2BBnX2NVPpfYcbzekHd9QW6fKNP2WRkRS5ZdWQIEFf3ZUnStsu6dv2%2BMK3KDCafKsNwOAi3j%2BkuQqb9k%2FtY4pcu393J%2B6HmCMW3hVu75YL0VeXO1SGmInLNj7XWj%2FKMUMR%2F0%2F3en2CwTQ9jwxI4UynS7zwc%2BnVCj4IQdqZd%2BQl0b2bx0hWoZMFJ987fbYP2OBtDzDaFZDDzUx%2Bt3FK%2FWj54%2BMNlbDpaeU39LcChCam7zUurbX9UvcPFOtQNTUZVijmpELyX%2BO%2B1Bv%2By%2BI7r8TNJyT68OdFUyhEi6%2BwN1dYt%2BEn3kUEQPiPnejgy0Q7CZkstpivDe%2FGf4PF6k7R78gpqoXGIkOGNq%2F25nYsSeHQSq%2B%2BJkYA%2FjOUaT30GuQHGgfpMaScWK4s1IOu%2BCfoxdeoh2aErPR7CmqzxtXeTSg2ajeqk%2BgBvBgAf0Q8JDP5UN0N3hZAVK8J1NujBil3Tu0HtVygYlJ0ActkSTfNN8qYgsr0jwzUEzCogL9ialOlFnpqX40SBNJkvVnQiNcuGXL06cLX6Q7JzngET2KBUfi%2B7o7p6IzDgMaTY89UGUFJfs1awnI41IwbOuMK3Spb8GOqQBnzEJTGRTzigBu4s%2B6iXfig3vLNPtqJjKfKXMhE3r%2Fl4rIrp0t58kwa1o%2BAd4F8K2HRvtkpeKGp2InFsYLguBd3YW9mI8n2xud1elLbWRyBe0wA6xZXXVcIxVz5WEl6gSJQ5114EntESUIPAjo9vaqKayHR0rjXxus9CWvF1eCj9kgryAw%2F0EpUz3h4PsWuGya%2F6helnJCzqXrQTnkN0CldpsFEo%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2003%3A00%3A11%20GMT&X-Amz-Signature=966e1638f9eb67c5f451aa7d989eed9eccacd96bc3271092b20c99e47b4d4a9b'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSI5RB6ZGB%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T190011Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJIMEYCIQCw9VcB8s63bA3LOPNHjnRYlPuYo07C22Nhbw7CzBsyogIhAMo53igAGwPughrSLjn9PT2aC%2Bn8cieuF2ARYHGRU0Q%2BKuYDCJH%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgyT3Etg0ZfLnQA5s%2BAqugPb9%2BBnX2NVPpfYcbzekHd9QW6fKNP2WRkRS5ZdWQIEFf3ZUnStsu6dv2%2BMK3KDCafKsNwOAi3j%2BkuQqb9k%2FtY4pcu393J%2B6HmCMW3hVu75YL0VeXO1SGmInLNj7XWj%2FKMUMR%2F0%2F3en2CwTQ9jwxI4UynS7zwc%2BnVCj4IQdqZd%2BQl0b2bx0hWoZMFJ987fbYP2OBtDzDaFZDDzUx%2Bt3FK%2FWj54%2BMNlbDpaeU39LcChCam7zUurbX9UvcPFOtQNTUZVijmpELyX%2BO%2B1Bv%2By%2BI7r8TNJyT68OdFUyhEi6%2BwN1dYt%2BEn3kUEQPiPnejgy0Q7CZkstpivDe%2FGf4PF6k7R78gpqoXGIkOGNq%2F25nYsSeHQSq%2B%2BJkYA%2FjOUaT30GuQHGgfpMaScWK4s1IOu%2BCfoxdeoh2aErPR7CmqzxtXeTSg2ajeqk%2BgBvBgAf0Q8JDP5UN0N3hZAVK8J1NujBil3Tu0HtVygYlJ0ActkSTfNN8qYgsr0jwzUEzCogL9ialOlFnpqX40SBNJkvVnQiNcuGXL06cLX6Q7JzngET2KBUfi%2B7o7p6IzDgMaTY89UGUFJfs1awnI41IwbOuMK3Spb8GOqQBnzEJTGRTzigBu4s%2B6iXfig3vLNPtqJjKfKXMhE3r%2Fl4rIrp0t58kwa1o%2BAd4F8K2HRvtkpeKGp2InFsYLguBd3YW9mI8n2xud1elLbWRyBe0wA6xZXXVcIxVz5WEl6gSJQ5114EntESUIPAjo9vaqKayHR0rjXxus9CWvF1eCj9kgryAw%2F0EpUz3h4PsWuGya%2F6helnJCzqXrQTnkN0CldpsFEo%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2003%3A00%3A11%20GMT&X-Amz-Signature=19dec86b9ff1bae81a817dbf2520a3438b8d72e7502b2d706d277bf6a4807ba0'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSI5RB6ZGB%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T190011Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJIMEYCIQCw9VcB8s63bA3LOPNHjnRYlPuYo07C22Nhbw7CzBsyogIhAMo53igAGwPughrSLjn9PT2aC%2Bn8cieuF2ARYHGRU0Q%2BKuYDCJH%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgyT3Etg0ZfLnQA5s%2BAqugPb9%2BBnX2NVPpfYcbzekHd9QW6fKNP2WRkRS5ZdWQIEFf3ZUnStsu6dv2%2BMK3KDCafKsNwOAi3j%2BkuQqb9k%2FtY4pcu393J%2B6HmCMW3hVu75YL0VeXO1SGmInLNj7XWj%2FKMUMR%2F0%2F3en2CwTQ9jwxI4UynS7zwc%2BnVCj4IQdqZd%2BQl0b2bx0hWoZMFJ987fbYP2OBtDzDaFZDDzUx%2Bt3FK%2FWj54%2BMNlbDpaeU39LcChCam7zUurbX9UvcPFOtQNTUZVijmpELyX%2BO%2B1Bv%2By%2BI7r8TNJyT68OdFUyhEi6%2BwN1dYt%2BEn3kUEQPiPnejgy0Q7CZkstpivDe%2FGf4PF6k7R78gpqoXGIkOGNq%2F25nYsSeHQSq%2B%2BJkYA%2FjOUaT30GuQHGgfpMaScWK4s1IOu%2BCfoxdeoh2aErPR7CmqzxtXeTSg2ajeqk%2BgBvBgAf0Q8JDP5UN0N3hZAVK8J1NujBil3Tu0HtVygYlJ0ActkSTfNN8qYgsr0jwzUEzCogL9ialOlFnpqX40SBNJkvVnQiNcuGXL06cLX6Q7JzngET2KBUfi%2B7o7p6IzDgMaTY89UGUFJfs1awnI41IwbOuMK3Spb8GOqQBnzEJTGRTzigBu4s%2B6iXfig3vLNPtqJjKfKXMhE3r%2Fl4rIrp0t58kwa1o%2BAd4F8K2HRvtkpeKGp2InFsYLguBd3YW9mI8n2xud1elLbWRyBe0wA6xZXXVcIxVz5WEl6gSJQ5114EntESUIPAjo9vaqKayHR0rjXxus9CWvF1eCj9kgryAw%2F0EpUz3h4PsWuGya%2F6helnJCzqXrQTnkN0CldpsFEo%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2003%3A00%3A11%20GMT&X-Amz-Signature=1139cfcf43649b1e8cedec32648a3df35ca23297f0678dea56cc27b2e5da3679'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSI5RB6ZGB%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T190011Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJIMEYCIQCw9VcB8s63bA3LOPNHjnRYlPuYo07C22Nhbw7CzBsyogIhAMo53igAGwPughrSLjn9PT2aC%2Bn8cieuF2ARYHGRU0Q%2BKuYDCJH%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgyT3Etg0ZfLnQA5s%2BAqugPb9%2BBnX2NVPpfYcbzekHd9QW6fKNP2WRkRS5ZdWQIEFf3ZUnStsu6dv2%2BMK3KDCafKsNwOAi3j%2BkuQqb9k%2FtY4pcu393J%2B6HmCMW3hVu75YL0VeXO1SGmInLNj7XWj%2FKMUMR%2F0%2F3en2CwTQ9jwxI4UynS7zwc%2BnVCj4IQdqZd%2BQl0b2bx0hWoZMFJ987fbYP2OBtDzDaFZDDzUx%2Bt3FK%2FWj54%2BMNlbDpaeU39LcChCam7zUurbX9UvcPFOtQNTUZVijmpELyX%2BO%2B1Bv%2By%2BI7r8TNJyT68OdFUyhEi6%2BwN1dYt%2BEn3kUEQPiPnejgy0Q7CZkstpivDe%2FGf4PF6k7R78gpqoXGIkOGNq%2F25nYsSeHQSq%2B%2BJkYA%2FjOUaT30GuQHGgfpMaScWK4s1IOu%2BCfoxdeoh2aErPR7CmqzxtXeTSg2ajeqk%2BgBvBgAf0Q8JDP5UN0N3hZAVK8J1NujBil3Tu0HtVygYlJ0ActkSTfNN8qYgsr0jwzUEzCogL9ialOlFnpqX40SBNJkvVnQiNcuGXL06cLX6Q7JzngET2KBUfi%2B7o7p6IzDgMaTY89UGUFJfs1awnI41IwbOuMK3Spb8GOqQBnzEJTGRTzigBu4s%2B6iXfig3vLNPtqJjKfKXMhE3r%2Fl4rIrp0t58kwa1o%2BAd4F8K2HRvtkpeKGp2InFsYLguBd3YW9mI8n2xud1elLbWRyBe0wA6xZXXVcIxVz5WEl6gSJQ5114EntESUIPAjo9vaqKayHR0rjXxus9CWvF1eCj9kgryAw%2F0EpUz3h4PsWuGya%2F6helnJCzqXrQTnkN0CldpsFEo%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2003%3A00%3A11%20GMT&X-Amz-Signature=9d86ae5bbfa1b7eeb3096a04a6fe49f005237f273421ff8c8f83a4dc7cb4312f'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSI5RB6ZGB%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T190011Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJIMEYCIQCw9VcB8s63bA3LOPNHjnRYlPuYo07C22Nhbw7CzBsyogIhAMo53igAGwPughrSLjn9PT2aC%2Bn8cieuF2ARYHGRU0Q%2BKuYDCJH%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgyT3Etg0ZfLnQA5s%2BAqugPb9%2BBnX2NVPpfYcbzekHd9QW6fKNP2WRkRS5ZdWQIEFf3ZUnStsu6dv2%2BMK3KDCafKsNwOAi3j%2BkuQqb9k%2FtY4pcu393J%2B6HmCMW3hVu75YL0VeXO1SGmInLNj7XWj%2FKMUMR%2F0%2F3en2CwTQ9jwxI4UynS7zwc%2BnVCj4IQdqZd%2BQl0b2bx0hWoZMFJ987fbYP2OBtDzDaFZDDzUx%2Bt3FK%2FWj54%2BMNlbDpaeU39LcChCam7zUurbX9UvcPFOtQNTUZVijmpELyX%2BO%2B1Bv%2By%2BI7r8TNJyT68OdFUyhEi6%2BwN1dYt%2BEn3kUEQPiPnejgy0Q7CZkstpivDe%2FGf4PF6k7R78gpqoXGIkOGNq%2F25nYsSeHQSq%2B%2BJkYA%2FjOUaT30GuQHGgfpMaScWK4s1IOu%2BCfoxdeoh2aErPR7CmqzxtXeTSg2ajeqk%2BgBvBgAf0Q8JDP5UN0N3hZAVK8J1NujBil3Tu0HtVygYlJ0ActkSTfNN8qYgsr0jwzUEzCogL9ialOlFnpqX40SBNJkvVnQiNcuGXL06cLX6Q7JzngET2KBUfi%2B7o7p6IzDgMaTY89UGUFJfs1awnI41IwbOuMK3Spb8GOqQBnzEJTGRTzigBu4s%2B6iXfig3vLNPtqJjKfKXMhE3r%2Fl4rIrp0t58kwa1o%2BAd4F8K2HRvtkpeKGp2InFsYLguBd3YW9mI8n2xud1elLbWRyBe0wA6xZXXVcIxVz5WEl6gSJQ5114EntESUIPAjo9vaqKayHR0rjXxus9CWvF1eCj9kgryAw%2F0EpUz3h4PsWuGya%2F6helnJCzqXrQTnkN0CldpsFEo%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2003%3A00%3A11%20GMT&X-Amz-Signature=80c6cde957637854ffab91aa03dfb41ba6678aa582d4811993a6e7adc740d5f2'
